# gray-edges.pptx: slide 1 gets an explicit accent1-colored background
# (previously it just inherited the master/theme background).
#
# Equivalent, in the Slide Background dialog, to:
#   Format Background -> Solid fill -> Color -> Theme Colors -> Accent 1
#   (and NOT clicking "Apply to All", so only slide 1 is affected).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Stop inheriting the layout/master background so this slide carries its
# own <p:bg> override.
$s.FollowMasterBackground = $false

# Make sure the background fill is a solid fill, then point it at the
# "Accent 1" theme color (msoThemeColorAccent1 = 5).
$s.Background.Fill.Solid()
$s.Background.Fill.ForeColor.ObjectThemeColor = 5
